# Automatische test-sync: 2025-06-22 19:05:50
#
# Adds the new incoming mail-log entry (row 35) to the "Logs" sheet and
# refreshes the "Dashboard" category breakdown to reflect the new counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 35
$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie openingstijden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Wanneer zijn jullie geopend?"
$logs.Cells.Item($newRow, 4).Value = "Openingstijden / Locatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor uw vraag. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Op zaterdag zijn wij gesloten. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 19:05:12"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# Extend the conditional-formatting ranges (Categorie / Beantwoord columns)
# so the new row inherits the same highlighting rules as the rest of the
# table (was D2:D34 / G2:G34, now D2:D35 / G2:G35).
$catRules = $logs.Range("D2:D34").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D35"))
}

$answeredRules = $logs.Range("G2:G34").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G35"))
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: the new entry's category ("Openingstijden / Locatie")
#    now ties at count 2, so it re-sorts above "Sollicitatie / Vacature"
#    (whose own count climbs to 2 as well); the two rows below shift down.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(9, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(9, 2).Value = 2

$dash.Cells.Item(10, 1).Value = "Sollicitatie / Vacature"
$dash.Cells.Item(10, 2).Value = 2

$dash.Cells.Item(11, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(11, 2).Value = 1

$dash.Cells.Item(12, 1).Value = "Uitnodiging / Evenement"
$dash.Cells.Item(12, 2).Value = 1
